$d = $word.ActiveDocument
$d.Content.Find.Execute("CU17", $true, $false, $false, $false, $false,
                         $true, 1, $false, "CU16", 2)
